$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.807.66'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '3.457.72'
$ws.Range("E3").Value = '  +2.23%  '
$ws.Range("E4").Value = '  -0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '582.95'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.48%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '147.51'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +7.31%  '
$ws.Range("D7").Value = '3.458.50'
$ws.Range("E7").Value = '  +2.31%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +1.45%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.68'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("E11").Value = '  +3.47%  '
$ws.Range("E12").Value = '  +2.75%  '
$ws.Range("D13").Value = '4.048.84'
$ws.Range("E13").Value = '  +2.28%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '27.99'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +9.09%  '
$ws.Range("E15").Value = '  -0.89%  '
$ws.Range("E16").Value = '  +1.54%  '
$ws.Range("D17").Value = '3.458.25'
$ws.Range("E17").Value = '  +2.28%  '
$ws.Range("D18").Value = '61.897.31'
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("E19").Value = '  +8.59%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.41'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +4.14%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '9.58'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +3.02%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '389.65'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +3.39%  '
$ws.Range("E23").Value = '  +2.82%  '
$ws.Range("E24").Value = '  +3.79%  '
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  -2.31%  '
$ws.Range("D28").Value = '3.601.61'
$ws.Range("E28").Value = '  +2.39%  '
$ws.Range("E29").Value = '  +0.97%  '
$ws.Range("E30").Value = '  +3.93%  '
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("E32").Value = '  -10.59%  '
$ws.Range("E33").Value = '  +1.71%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +2.23%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '24.23'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("D37").Value = '3.486.86'
$ws.Range("E37").Value = '  +2.32%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '7.02'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +2.86%  '
$ws.Range("E39").Value = '  +1.70%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '5.16'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.41%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '167.05'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.41%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0785'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +3.61%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '27.41'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +6.32%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.804'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +3.95%  '
$ws.Range("E45").Value = '  +2.27%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '4.51'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +3.78%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("D50").Value = '2.576.25'
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("E51").Value = '  +2.47%  '
